$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.272.66'
$ws.Range('E2').Value = '  +8.35%  '
$ws.Range('D3').Value = '1.598.37'
$ws.Range('E3').Value = '  +8.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9966'
$ws.Range('E5').Value = '  +2.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '288.82'
$ws.Range('E6').Value = '  +3.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3711'
$ws.Range('E7').Value = '  +1.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3376'
$ws.Range('E8').Value = '  +9.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.33'
$ws.Range('E9').Value = '  +6.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.140'
$ws.Range('E10').Value = '  +7.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07047'
$ws.Range('E11').Value = '  +5.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.79'
$ws.Range('E13').Value = '  +9.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.920'
$ws.Range('E14').Value = '  +7.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.657'
$ws.Range('E15').Value = '  +7.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001082'
$ws.Range('E16').Value = '  +4.98%  '
$ws.Range('B17').Value = 'Dai'
$ws.Range('C17').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9959'
$ws.Range('E17').Value = '  +2.19%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '1.599.21'
$ws.Range('E18').Value = '  +8.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06636'
$ws.Range('E19').Value = '  +11.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '78.09'
$ws.Range('E20').Value = '  +12.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.18'
$ws.Range('E21').Value = '  +11.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.047'
$ws.Range('E22').Value = '  +10.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.77'
$ws.Range('E23').Value = '  +6.68%  '
$ws.Range('D24').Value = '22.352.17'
$ws.Range('E24').Value = '  +8.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.394'
$ws.Range('E25').Value = '  +6.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.490'
$ws.Range('E26').Value = '  +16.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.14'
$ws.Range('E27').Value = '  +7.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.69'
$ws.Range('E28').Value = '  +14.24%  '
$ws.Range('D29').Value = '1.781.30'
$ws.Range('E29').Value = '  +9.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.47'
$ws.Range('E30').Value = '  +5.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.190'
$ws.Range('E31').Value = '  +5.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.967'
$ws.Range('E32').Value = '  +19.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9455'
$ws.Range('E33').Value = '  +16.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08258'
$ws.Range('E34').Value = '  +2.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.612'
$ws.Range('E35').Value = '  +4.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.300'
$ws.Range('E36').Value = '  +12.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.647'
$ws.Range('E37').Value = '  +13.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '11.77'
$ws.Range('E38').Value = '  +13.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06146'
$ws.Range('E39').Value = '  +5.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.240'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02213'
$ws.Range('E41').Value = '  +8.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2022'
$ws.Range('E42').Value = '  +7.32%  '
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9952'
$ws.Range('E43').Value = '  +2.14%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5946'
$ws.Range('E44').Value = '  +12.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.15'
$ws.Range('E45').Value = '  +8.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.662'
$ws.Range('E46').Value = '  +4.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5719'
$ws.Range('E47').Value = '  +10.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.01'
$ws.Range('E48').Value = '  +4.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.974'
$ws.Range('E49').Value = '  +10.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06821'
$ws.Range('E50').Value = '  +5.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.69'
$ws.Range('E51').Value = '  +9.35%  '
